$wb = $excel.ActiveWorkbook

# Both "展览" (Exhibition) and "全部类型" (All Types) sheets contain the same
# event rows and both need the "想去人数" (number of people interested)
# values updated for row 5 and row 6.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 838
    $ws.Range("F6").Value = 209
}
